$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: each entry is Row number + new values for B (Coin), C (Link), D (Price), E (Volume)
$updates = @(
    @{Row=2; D='30.698.04'; E='  +1.89%  '},
    @{Row=3; D='1.900.48'; E='  +2.86%  '},
    @{Row=4; D='0.9995'; E='  -0.15%  '},
    @{Row=5; D='239.22'; E='  +1.47%  '},
    @{Row=6; D='1.000'; E='  -0.09%  '},
    @{Row=7; D='0.4817'; E='  +1.37%  '},
    @{Row=8; D='0.2844'; E='  +1.00%  '},
    @{Row=9; D='0.06566'; E='  +1.49%  '},
    @{Row=10; D='1.923.04'; E='  +3.75%  '},
    @{Row=11; D='0.07471'; E='  +2.42%  '},
    @{Row=12; D='16.73'; E='  +2.67%  '},
    @{Row=13; D='5.124'; E='  +0.35%  '},
    @{Row=14; D='88.23'; E='  +1.34%  '},
    @{Row=15; D='0.6679'; E='  +3.64%  '},
    @{Row=16; D='30.668.05'; E='  +1.97%  '},
    @{Row=17; D='13.34'; E='  +0.99%  '},
    @{Row=18; D='1.000'; E='  -0.07%  '},
    @{Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007626'; E='  +0.31%  '},
    @{Row=20; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='233.51'; E='  +8.13%  '},
    @{Row=21; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.186.27'; E='  +3.51%  '},
    @{Row=22; D='5.306'; E='  +0.97%  '},
    @{Row=23; D='0.9997'; E='  -0.18%  '},
    @{Row=24; D='6.231'; E='  +2.27%  '},
    @{Row=25; D='170.37'; E='  +4.58%  '},
    @{Row=26; D='9.309'; E='  +1.42%  '},
    @{Row=27; D='18.68'; E='  +2.01%  '},
    @{Row=28; D='1.965'; E='  +2.76%  '},
    @{Row=29; E='  -1.93%  '},
    @{Row=30; E='  +9.59%  '},
    @{Row=31; D='4.376'; E='  +3.56%  '},
    @{Row=32; D='4.042'; E='  +2.13%  '},
    @{Row=33; D='0.05085'; E='  +1.56%  '},
    @{Row=34; D='1.218'; E='  +7.69%  '},
    @{Row=35; D='0.7556'; E='  +2.29%  '},
    @{Row=36; E='  +1.07%  '},
    @{Row=37; D='0.01886'; E='  +3.80%  '},
    @{Row=38; D='2.658'; E='  +2.21%  '},
    @{Row=39; D='0.9183'; E='  +2.08%  '},
    @{Row=40; D='2.086'; E='  +1.38%  '},
    @{Row=41; D='107.24'; E='  +1.11%  '},
    @{Row=42; D='0.4312'; E='  +1.80%  '},
    @{Row=43; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.005'; E='  +0.57%  '},
    @{Row=44; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='5.786'; E='  -2.04%  '},
    @{Row=45; D='7.442'; E='  +0.45%  '},
    @{Row=46; D='64.75'; E='  +1.50%  '},
    @{Row=47; D='0.1278'; E='  -1.94%  '},
    @{Row=48; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.486'; E='  -3.99%  '},
    @{Row=49; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.027'; E='  +3.50%  '},
    @{Row=50; D='33.96'; E='  -0.52%  '},
    @{Row=51; D='0.05671'; E='  -0.12%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        $c = $ws.Range("B" + $r)
        $c.Value = "'" + $u.B
        $c.Style = "Normal"
    }
    if ($u.ContainsKey("C")) {
        $c = $ws.Range("C" + $r)
        $c.Value = "'" + $u.C
        $c.Style = "Normal"
    }
    if ($u.ContainsKey("D")) {
        $c = $ws.Range("D" + $r)
        $c.Value = "'" + $u.D
        $c.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $c = $ws.Range("E" + $r)
        $c.Value = $u.E
    }
}

Write-Output "Updated $($updates.Count) rows"
